$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.544.39'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '1.564.81'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.68'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.494'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.34'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +5.30%  '
$ws.Range("E9").Value = '  +0.91%  '
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0883'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("E13").Value = '  -1.72%  '
$ws.Range("D14").Value = '1.566.18'
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("E15").Value = '  -1.97%  '
$ws.Range("D16").Value = '28.546.04'
$ws.Range("E16").Value = '  +0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.67'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '62.10'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.76'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.07%  '
$ws.Range("D20").Value = '0.0₃0694'
$ws.Range("E20").Value = '  -2.25%  '
$ws.Range("E21").Value = '  -2.55%  '
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.86'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -6.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.14'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.06'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +5.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.34'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.95'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.30%  '
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("E29").Value = '  -3.74%  '
$ws.Range("E31").Value = '  -1.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.10'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.00%  '
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("D35").Value = '1.396.88'
$ws.Range("E35").Value = '  -1.20%  '
$ws.Range("E36").Value = '  -1.28%  '
$ws.Range("E37").Value = '  -3.79%  '
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.59'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.26%  '
$ws.Range("E40").Value = '  -1.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.536'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.31%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("E43").Value = '  -3.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.87'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.61%  '
$ws.Range("E45").Value = '  -4.74%  '
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.53'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -3.41%  '
$ws.Range("D48").Value = '1.701.50'
$ws.Range("E48").Value = '  -1.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.20'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("E50").Value = '  -3.94%  '
$ws.Range("E51").Value = '  -0.96%  '
